# Auto-generated Excel COM-interop script
# Applies numeric cell updates to multiple sheets per the target diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 261.85715
$ws.Range("I2").Value = 96
$ws.Range("J2").Value = 307.0909
$ws.Range("K2").Value = 96
$ws.Range("L2").Value = 307.0909
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = -533.0908999999999
$ws.Range("H33").Value = 617.6
$ws.Range("I33").Value = 632.8889
$ws.Range("K33").Value = 632.8889
$ws.Range("M33").Value = -403.8889
$ws.Range("H100").Value = 55557310
$ws.Range("I100").Value = 83334984
$ws.Range("J100").Value = 1962.6666
$ws.Range("K100").Value = 83334984
$ws.Range("L100").Value = 1962.6666
$ws.Range("M100").Value = -83334443
$ws.Range("N100").Value = -3044.6666
$ws.Range("H112").Value = 2151.75
$ws.Range("J112").Value = 2445.182
$ws.Range("L112").Value = 7335.545999999999
$ws.Range("N112").Value = -9551.545999999998
$ws.Range("H129").Value = 866.36
$ws.Range("I129").Value = 312.8889
$ws.Range("J129").Value = 921.0989
$ws.Range("K129").Value = 938.6667
$ws.Range("L129").Value = 2763.2967
$ws.Range("M129").Value = 4061.3333
$ws.Range("N129").Value = -12763.2967
$ws.Range("H137").Value = 1816.4286
$ws.Range("I137").Value = 1002
$ws.Range("J137").Value = 2038.5454
$ws.Range("K137").Value = 3006
$ws.Range("L137").Value = 6115.6362
$ws.Range("M137").Value = -456
$ws.Range("N137").Value = -11215.6362
$ws.Range("H138").Value = 1501.9192
$ws.Range("J138").Value = 2006.4839
$ws.Range("L138").Value = 6019.4517
$ws.Range("N138").Value = -16299.4517
$ws.Range("J141").Value = 795
$ws.Range("L141").Value = 2385
$ws.Range("N141").Value = -12745

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2843.3206
$ws.Range("I32").Value = 2643.9077
$ws.Range("K32").Value = 2643.9077
$ws.Range("M32").Value = -2356.9077
$ws.Range("H61").Value = 1045.4722
$ws.Range("I61").Value = 989.2857
$ws.Range("J61").Value = 1242.125
$ws.Range("K61").Value = 989.2857
$ws.Range("L61").Value = 1242.125
$ws.Range("M61").Value = -777.2857
$ws.Range("N61").Value = -1666.125
$ws.Range("H64").Value = 59999.332
$ws.Range("J64").Value = 59999.332
$ws.Range("L64").Value = 59999.332
$ws.Range("N64").Value = -60495.332
$ws.Range("H67").Value = 59999.332
$ws.Range("J67").Value = 59999.332
$ws.Range("L67").Value = 59999.332
$ws.Range("N67").Value = -61715.332
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H114").Value = 23316.084
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 23316.084
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 23316.084
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -31994.084
$ws.Range("H132").Value = 1485.8788
$ws.Range("I132").Value = 1144.1428
$ws.Range("K132").Value = 3432.4284
$ws.Range("M132").Value = -902.4284000000002
$ws.Range("H136").Value = 1045.4722
$ws.Range("I136").Value = 989.2857
$ws.Range("J136").Value = 1242.125
$ws.Range("K136").Value = 2967.8571
$ws.Range("L136").Value = 3726.375
$ws.Range("M136").Value = -417.8571000000002
$ws.Range("N136").Value = -8826.375

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1192.5283
$ws.Range("I31").Value = 1174.2
$ws.Range("J31").Value = 1498
$ws.Range("K31").Value = 1174.2
$ws.Range("L31").Value = 1498
$ws.Range("M31").Value = -879.2
$ws.Range("N31").Value = -2088
$ws.Range("H34").Value = 1192.5283
$ws.Range("I34").Value = 1174.2
$ws.Range("J34").Value = 1498
$ws.Range("K34").Value = 1174.2
$ws.Range("L34").Value = 1498
$ws.Range("M34").Value = -972.2
$ws.Range("N34").Value = -1902
$ws.Range("H122").Value = 1039.9231
$ws.Range("I122").Value = 1040
$ws.Range("J122").Value = 1039.8
$ws.Range("K122").Value = 3120
$ws.Range("L122").Value = 3119.4
$ws.Range("M122").Value = -670
$ws.Range("N122").Value = -8019.4
$ws.Range("H132").Value = 1889.05
$ws.Range("I132").Value = 984.5
$ws.Range("K132").Value = 2953.5
$ws.Range("M132").Value = -423.5
$ws.Range("H134").Value = 1312.8
$ws.Range("I134").Value = 1013.55554
$ws.Range("J134").Value = 1761.6666
$ws.Range("K134").Value = 3040.66662
$ws.Range("L134").Value = 5284.9998
$ws.Range("M134").Value = -505.66662
$ws.Range("N134").Value = -10354.9998

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5250
$ws.Range("J80").Value = 5250
$ws.Range("L80").Value = 15750
$ws.Range("N80").Value = -17622
$ws.Range("H83").Value = 5250
$ws.Range("J83").Value = 5250
$ws.Range("L83").Value = 47250
$ws.Range("N83").Value = -56610
$ws.Range("H130").Value = 2032.8572
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("H131").Value = 12050337
$ws.Range("I131").Value = 111111540
$ws.Range("J131").Value = 2353.7163
$ws.Range("K131").Value = 333334620
$ws.Range("L131").Value = 7061.1489
$ws.Range("M131").Value = -333329580
$ws.Range("N131").Value = -17141.1489
$ws.Range("H137").Value = 2452.72
$ws.Range("I137").Value = 1162.5
$ws.Range("J137").Value = 3643.6924
$ws.Range("K137").Value = 3487.5
$ws.Range("L137").Value = 10931.0772
$ws.Range("M137").Value = 1612.5
$ws.Range("N137").Value = -21131.0772
$ws.Range("H139").Value = 1666.9259
$ws.Range("I139").Value = 1643.7826
$ws.Range("J139").Value = 1800
$ws.Range("K139").Value = 4931.3478
$ws.Range("L139").Value = 5400
$ws.Range("M139").Value = 208.6522000000004
$ws.Range("N139").Value = -15680

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 1806.3889
$ws.Range("I122").Value = 2026.1818
$ws.Range("J122").Value = 1461
$ws.Range("K122").Value = 6078.5454
$ws.Range("L122").Value = 4383
$ws.Range("M122").Value = -3628.5454
$ws.Range("N122").Value = -9283
$ws.Range("H132").Value = 2041.2927
$ws.Range("I132").Value = 1494.9565
$ws.Range("J132").Value = 2739.389
$ws.Range("K132").Value = 4484.8695
$ws.Range("L132").Value = 8218.167000000001
$ws.Range("M132").Value = -1954.8695
$ws.Range("N132").Value = -13278.167

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2152.4644
$ws.Range("I40").Value = 1709.6666
$ws.Range("J40").Value = 2949.5
$ws.Range("K40").Value = 1709.6666
$ws.Range("L40").Value = 2949.5
$ws.Range("M40").Value = -1573.6666
$ws.Range("N40").Value = -3221.5
$ws.Range("H100").Value = 1749
$ws.Range("I100").Value = 998
$ws.Range("K100").Value = 998
$ws.Range("M100").Value = -457

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("H20").Value = 70011
$ws.Range("J20").Value = 70011
$ws.Range("L20").Value = 70011
$ws.Range("N20").Value = -70491
$ws.Range("H119").Value = 17979.2
$ws.Range("J119").Value = 17979.2
$ws.Range("L119").Value = 17979.2
$ws.Range("N119").Value = -27655.2
